# Helper: force a (possibly numeric-looking) string into a cell as TEXT,
# stored as a shared string, without altering the cell style/number format.
# Plain numeric-looking strings assigned via .Value get auto-converted to
# real numbers by Excel, so for those we build a text formula and then
# convert the formula to a static value via copy / paste-special-values.
function Set-TextValue($range, [string]$val) {
    $range.Formula = '="' + $val + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

$wb = $excel.ActiveWorkbook

# --- Worksheet "Restricciones_del_lider" (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "-16.45 + x_1 + x_2 + y_1 - 2y_2"
Set-TextValue $ws.Range("B2") "-23.55"
Set-TextValue $ws.Range("D2") "0.86"

# --- Worksheet "Restricciones_del_follower" (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "18.85 - x_1 + 2y_2"
Set-TextValue $ws.Range("B2") "-8.850000000000001"
Set-TextValue $ws.Range("D2") "0.32"
Set-TextValue $ws.Range("E2") "0"
Set-TextValue $ws.Range("F2") "8.9"
$ws.Range("A3").Value = "20.0 - x_2 + 2y_2"
Set-TextValue $ws.Range("B3") "-10.0"
Set-TextValue $ws.Range("D3") "0.82"
Set-TextValue $ws.Range("F3") "6.4"
$ws.Range("A4").Value = "-5.1 - y_1"
Set-TextValue $ws.Range("B4") "-4.9"
Set-TextValue $ws.Range("D4") "0.26"
Set-TextValue $ws.Range("E4") "-2.9"
Set-TextValue $ws.Range("F4") "-8.6"
$ws.Range("A5").Value = "5.100000000000001 + y_1"
Set-TextValue $ws.Range("B5") "-25.1"
Set-TextValue $ws.Range("D5") "0.85"
Set-TextValue $ws.Range("E5") "0"
Set-TextValue $ws.Range("F5") "1.9"
$ws.Range("A6").Value = "-11.35 - y_2"
Set-TextValue $ws.Range("B6") "-1.3499999999999996"
Set-TextValue $ws.Range("D6") "0.7"
Set-TextValue $ws.Range("E6") "0"
Set-TextValue $ws.Range("F6") "5.6000000000000005"
$ws.Range("A7").Value = "-48.65 + y_2"
Set-TextValue $ws.Range("B7") "-28.65"
Set-TextValue $ws.Range("D7") "0.51"
Set-TextValue $ws.Range("E7") "7.9"
Set-TextValue $ws.Range("F7") "1.5"

# --- Worksheet "Punto_modificado" (index 4) ---
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Range("A2") "1.55"
Set-TextValue $ws.Range("B2") "2.7"
Set-TextValue $ws.Range("C2") "-5.1"
Set-TextValue $ws.Range("D2") "-8.65"

# --- Worksheet "Vector_bf" (index 5) ---
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "-27.29"
Set-TextValue $ws.Range("A3") "-19.39"

# --- Worksheet "Vector_BF" (index 6) ---
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "-2.86"
Set-TextValue $ws.Range("A3") "-2.86"
Set-TextValue $ws.Range("A4") "-0.7599999999999998"
Set-TextValue $ws.Range("A5") "-3.1800000000000006"
